$d = $word.ActiveDocument

$replacements = @(
    @{old = "2024-09-04 Wednesday"; new = "2024-09-05 Thursday"},
    @{old = "57×96="; new = "43×63="},
    @{old = "11×39="; new = "14×32="},
    @{old = "59×74="; new = "54×21="},
    @{old = "31×90="; new = "50×48="},
    @{old = "34×94="; new = "62×69="},
    @{old = "35×94="; new = "12×82="},
    @{old = "81×54="; new = "73×56="},
    @{old = "19×47="; new = "46×53="},
    @{old = "73×18="; new = "50×93="},
    @{old = "68×42="; new = "34×16="},
    @{old = "32×50="; new = "42×51="},
    @{old = "12×31="; new = "12×17="},
    @{old = "91×74="; new = "80×53="},
    @{old = "57×99="; new = "78×69="},
    @{old = "45×46="; new = "45×67="},
    @{old = "87×62="; new = "20×35="},
    @{old = "39×74="; new = "29×48="},
    @{old = "43×89="; new = "41×46="},
    @{old = "83×96="; new = "75×13="},
    @{old = "68×56="; new = "97×56="},
    @{old = "14×41="; new = "69×82="},
    @{old = "58×19="; new = "33×29="},
    @{old = "41×30="; new = "40×72="},
    @{old = "43×15="; new = "95×30="},
    @{old = "16×51="; new = "58×30="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

$d.Save()
